$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "special offer" value for the "Ркацители" wine row,
# matching the same promo text already used on rows 2 and 9.
$ws.Range("F5").Value = "Выгодное предложение"

# Update the selected cell as left by the author after editing.
$ws.Range("E15").Select()
